# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to Leve profit columns (H-N)
# across all sheets, per the authoritative diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 5150  # H18
$ws.Cells.Item(18, 9).Value = 5150  # I18
$ws.Cells.Item(18, 11).Value = 5150  # K18
$ws.Cells.Item(18, 13).Value = -4866  # M18
$ws.Cells.Item(40, 8).Value = 3000  # H40
$ws.Cells.Item(40, 9).Value = 0  # I40
$ws.Cells.Item(40, 11).Value = 0  # K40
$ws.Cells.Item(40, 13).Value = $null  # M40: clear (was -3658.25)
$ws.Cells.Item(51, 8).Value = 10420140  # H51
$ws.Cells.Item(51, 10).Value = 13892372  # J51
$ws.Cells.Item(51, 12).Value = 13892372  # L51
$ws.Cells.Item(51, 14).Value = -13893340  # N51
$ws.Cells.Item(64, 8).Value = 2977.7778  # H64
$ws.Cells.Item(67, 8).Value = 2977.7778  # H67
$ws.Cells.Item(112, 8).Value = 1996.8889  # H112
$ws.Cells.Item(112, 10).Value = 2078.739  # J112
$ws.Cells.Item(112, 12).Value = 6236.217000000001  # L112
$ws.Cells.Item(112, 14).Value = -8452.217000000001  # N112
$ws.Cells.Item(113, 8).Value = 81604.92  # H113
$ws.Cells.Item(113, 9).Value = 170310.83  # I113
$ws.Cells.Item(113, 11).Value = 170310.83  # K113
$ws.Cells.Item(113, 13).Value = -167056.83  # M113
$ws.Cells.Item(138, 8).Value = 3805.25  # H138
$ws.Cells.Item(138, 10).Value = 4010.8413  # J138
$ws.Cells.Item(138, 12).Value = 12032.5239  # L138
$ws.Cells.Item(138, 14).Value = -22312.5239  # N138

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2237.1  # H2
$ws.Cells.Item(2, 9).Value = 1819.5555  # I2
$ws.Cells.Item(2, 10).Value = 5995  # J2
$ws.Cells.Item(2, 11).Value = 1819.5555  # K2
$ws.Cells.Item(2, 12).Value = 5995  # L2
$ws.Cells.Item(2, 13).Value = -1706.5555  # M2
$ws.Cells.Item(2, 14).Value = -6221  # N2
$ws.Cells.Item(32, 8).Value = 8722.514999999999  # H32
$ws.Cells.Item(32, 9).Value = 6705.279  # I32
$ws.Cells.Item(32, 10).Value = 26301.285  # J32
$ws.Cells.Item(32, 11).Value = 6705.279  # K32
$ws.Cells.Item(32, 12).Value = 26301.285  # L32
$ws.Cells.Item(32, 13).Value = -6418.279  # M32
$ws.Cells.Item(32, 14).Value = -26875.285  # N32
$ws.Cells.Item(45, 8).Value = 5580.407  # H45
$ws.Cells.Item(45, 9).Value = 5988.143  # I45
$ws.Cells.Item(45, 11).Value = 5988.143  # K45
$ws.Cells.Item(45, 13).Value = -5611.143  # M45
$ws.Cells.Item(74, 8).Value = 1736.027  # H74
$ws.Cells.Item(74, 9).Value = 1428.7587  # I74
$ws.Cells.Item(74, 11).Value = 1428.7587  # K74
$ws.Cells.Item(74, 13).Value = -554.7587000000001  # M74
$ws.Cells.Item(77, 8).Value = 1736.027  # H77
$ws.Cells.Item(77, 9).Value = 1428.7587  # I77
$ws.Cells.Item(77, 11).Value = 7143.793500000001  # K77
$ws.Cells.Item(77, 13).Value = -2775.793500000001  # M77
$ws.Cells.Item(97, 8).Value = 791.86365  # H97
$ws.Cells.Item(97, 9).Value = 907.1177  # I97
$ws.Cells.Item(97, 11).Value = 907.1177  # K97
$ws.Cells.Item(97, 13).Value = -411.1177  # M97
$ws.Cells.Item(116, 8).Value = 2237.1  # H116
$ws.Cells.Item(116, 9).Value = 1819.5555  # I116
$ws.Cells.Item(116, 10).Value = 5995  # J116
$ws.Cells.Item(116, 11).Value = 1819.5555  # K116
$ws.Cells.Item(116, 12).Value = 5995  # L116
$ws.Cells.Item(116, 13).Value = 474.4445000000001  # M116
$ws.Cells.Item(116, 14).Value = -10583  # N116
$ws.Cells.Item(132, 8).Value = 5536.732  # H132
$ws.Cells.Item(132, 9).Value = 5860.9795  # I132
$ws.Cells.Item(132, 11).Value = 17582.9385  # K132
$ws.Cells.Item(132, 13).Value = -15052.9385  # M132

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2237.1  # H3
$ws.Cells.Item(3, 9).Value = 1819.5555  # I3
$ws.Cells.Item(3, 10).Value = 5995  # J3
$ws.Cells.Item(3, 11).Value = 1819.5555  # K3
$ws.Cells.Item(3, 12).Value = 5995  # L3
$ws.Cells.Item(3, 13).Value = -1705.5555  # M3
$ws.Cells.Item(3, 14).Value = -6223  # N3
$ws.Cells.Item(20, 8).Value = 15158235  # H20
$ws.Cells.Item(20, 9).Value = 20841938  # I20
$ws.Cells.Item(20, 11).Value = 20841938  # K20
$ws.Cells.Item(20, 13).Value = -20841691  # M20

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3808.1904  # H31
$ws.Cells.Item(31, 9).Value = 2807.6667  # I31
$ws.Cells.Item(31, 10).Value = 4908.7666  # J31
$ws.Cells.Item(31, 11).Value = 2807.6667  # K31
$ws.Cells.Item(31, 12).Value = 4908.7666  # L31
$ws.Cells.Item(31, 13).Value = -2512.6667  # M31
$ws.Cells.Item(31, 14).Value = -5498.7666  # N31
$ws.Cells.Item(34, 8).Value = 3808.1904  # H34
$ws.Cells.Item(34, 9).Value = 2807.6667  # I34
$ws.Cells.Item(34, 10).Value = 4908.7666  # J34
$ws.Cells.Item(34, 11).Value = 2807.6667  # K34
$ws.Cells.Item(34, 12).Value = 4908.7666  # L34
$ws.Cells.Item(34, 13).Value = -2605.6667  # M34
$ws.Cells.Item(34, 14).Value = -5312.7666  # N34
$ws.Cells.Item(58, 8).Value = 1720.1875  # H58
$ws.Cells.Item(58, 10).Value = 1800  # J58
$ws.Cells.Item(58, 12).Value = 1800  # L58
$ws.Cells.Item(58, 14).Value = -2206  # N58
$ws.Cells.Item(59, 8).Value = 17500  # H59
$ws.Cells.Item(60, 8).Value = 13650  # H60
$ws.Cells.Item(68, 8).Value = 18424.143  # H68
$ws.Cells.Item(68, 10).Value = 18424.143  # J68
$ws.Cells.Item(68, 12).Value = 18424.143  # L68
$ws.Cells.Item(68, 14).Value = -19922.143  # N68
$ws.Cells.Item(71, 8).Value = 18424.143  # H71
$ws.Cells.Item(71, 10).Value = 18424.143  # J71
$ws.Cells.Item(71, 12).Value = 55272.429  # L71
$ws.Cells.Item(71, 14).Value = -62760.429  # N71
$ws.Cells.Item(132, 8).Value = 2324.8206  # H132
$ws.Cells.Item(132, 9).Value = 1949.1936  # I132
$ws.Cells.Item(132, 11).Value = 5847.5808  # K132
$ws.Cells.Item(132, 13).Value = -3317.5808  # M132
$ws.Cells.Item(136, 8).Value = 1720.1875  # H136
$ws.Cells.Item(136, 10).Value = 1800  # J136
$ws.Cells.Item(136, 12).Value = 5400  # L136
$ws.Cells.Item(136, 14).Value = -10500  # N136

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(52, 8).Value = 0  # H52
$ws.Cells.Item(52, 9).Value = 0  # I52
$ws.Cells.Item(52, 10).Value = 0  # J52
$ws.Cells.Item(52, 11).Value = 0  # K52
$ws.Cells.Item(52, 12).Value = 0  # L52
$ws.Cells.Item(52, 13).Value = $null  # M52: clear (was -2999731)
$ws.Cells.Item(52, 14).Value = $null  # N52: clear (was -3502)
$ws.Cells.Item(56, 8).Value = 8857.950000000001  # H56
$ws.Cells.Item(56, 9).Value = 8857.950000000001  # I56
$ws.Cells.Item(56, 11).Value = 8857.950000000001  # K56
$ws.Cells.Item(56, 13).Value = -8327.950000000001  # M56
$ws.Cells.Item(87, 8).Value = 15103.5  # H87
$ws.Cells.Item(87, 9).Value = 5207  # I87
$ws.Cells.Item(87, 11).Value = 15621  # K87
$ws.Cells.Item(87, 13).Value = -14373  # M87
$ws.Cells.Item(90, 8).Value = 15103.5  # H90
$ws.Cells.Item(90, 9).Value = 5207  # I90
$ws.Cells.Item(90, 11).Value = 46863  # K90
$ws.Cells.Item(90, 13).Value = -40623  # M90

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 2747.25  # H113
$ws.Cells.Item(113, 9).Value = 2456.6  # I113
$ws.Cells.Item(113, 10).Value = 3231.6667  # J113
$ws.Cells.Item(113, 11).Value = 2456.6  # K113
$ws.Cells.Item(113, 12).Value = 3231.6667  # L113
$ws.Cells.Item(113, 13).Value = -286.5999999999999  # M113
$ws.Cells.Item(113, 14).Value = -7571.6667  # N113

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1813.3334  # H22
$ws.Cells.Item(22, 9).Value = 1865.6  # I22
$ws.Cells.Item(22, 10).Value = 1748  # J22
$ws.Cells.Item(22, 11).Value = 1865.6  # K22
$ws.Cells.Item(22, 12).Value = 1748  # L22
$ws.Cells.Item(22, 13).Value = -1570.6  # M22
$ws.Cells.Item(22, 14).Value = -2338  # N22
$ws.Cells.Item(27, 8).Value = 1813.3334  # H27
$ws.Cells.Item(27, 9).Value = 1865.6  # I27
$ws.Cells.Item(27, 10).Value = 1748  # J27
$ws.Cells.Item(27, 11).Value = 1865.6  # K27
$ws.Cells.Item(27, 12).Value = 1748  # L27
$ws.Cells.Item(27, 13).Value = -1758.6  # M27
$ws.Cells.Item(27, 14).Value = -1962  # N27
$ws.Cells.Item(46, 8).Value = 24065.1  # H46
$ws.Cells.Item(46, 9).Value = 40982.273  # I46
$ws.Cells.Item(46, 10).Value = 3388.5557  # J46
$ws.Cells.Item(46, 11).Value = 40982.273  # K46
$ws.Cells.Item(46, 12).Value = 3388.5557  # L46
$ws.Cells.Item(46, 13).Value = -40794.273  # M46
$ws.Cells.Item(46, 14).Value = -3764.5557  # N46
$ws.Cells.Item(93, 8).Value = 15057.414  # H93
$ws.Cells.Item(93, 9).Value = 3371.44  # I93
$ws.Cells.Item(93, 11).Value = 3371.44  # K93
$ws.Cells.Item(93, 13).Value = -2123.44  # M93

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 747630.3  # H132
$ws.Cells.Item(132, 9).Value = 1098673.9  # I132
$ws.Cells.Item(132, 11).Value = 3296021.7  # K132
$ws.Cells.Item(132, 13).Value = -3293491.7  # M132
